# Example_PlateMap.xlsx update
#
# Summary of the change (per the commit message / diff):
#   - "concentration" work continues: the Starting_Dilution_or_concentration
#     values for the second virus block on Sheet1 (rows 12-21, column J) are
#     re-prefilled from 20 -> 10 while validating plates with multiple
#     viruses.
#   - The last row of that block (J21) picks up the same "no bottom border"
#     formatting as the rest of the block (a side-effect of how the values
#     were filled down), while I21/K21 keep their original formatting.
#   - Selections / active sheet move around as the author works: Sheet1
#     becomes the active sheet/tab, with the cursor resting on J20; the
#     "concentration" sheet's selection moves to G43; "Sheet_with_errors"
#     is no longer the active tab.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("concentration")
$ws3 = $wb.Worksheets.Item("Sheet_with_errors")

# --- Sheet1: re-prefill the starting dilution/concentration column -------
# Rows 12-21 (the second virus block) go from 20 to the new default of 10.
$ws1.Range("J12:J21").Value = 10

# J21 picks up the interior-row formatting (no bottom border) that J12:J20
# already use, matching the look of the rest of the filled block - copy
# the format from J20 down onto J21 without touching J21's value.
$ws1.Range("J20").Copy()
$ws1.Range("J21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Selections / active sheet -------------------------------------------
# concentration sheet: selection moves to G43
$ws2.Activate()
$ws2.Range("G43").Select()

# Sheet_with_errors keeps its own selection (B2) but is no longer active.

# Sheet1 becomes the active sheet/tab, cursor resting on J20
$ws1.Activate()
$ws1.Range("J20").Select()
